$wb = $excel.ActiveWorkbook

# ALC row 6
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 121.85714
$ws.Range("I6").Value = 86.666664
$ws.Range("J6").Value = 333
$ws.Range("K6").Value = 259.999992
$ws.Range("L6").Value = 999
$ws.Range("M6").Value = -147.999992
$ws.Range("N6").Value = -1223

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1078.4117
$ws.Range("I33").Value = 1139.4375
$ws.Range("J33").Value = 102
$ws.Range("K33").Value = 1139.4375
$ws.Range("L33").Value = 102
$ws.Range("M33").Value = -910.4375
$ws.Range("N33").Value = -560

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1987.2858
$ws.Range("I100").Value = 1962.3334
$ws.Range("J100").Value = 2006
$ws.Range("K100").Value = 1962.3334
$ws.Range("L100").Value = 2006
$ws.Range("M100").Value = -1421.3334
$ws.Range("N100").Value = -3088

# ALC row 121
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 7253
$ws.Range("J121").Value = 6829.5
$ws.Range("L121").Value = 20488.5
$ws.Range("N121").Value = -23982.5

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 10205.174
$ws.Range("I137").Value = 11649.066
$ws.Range("J137").Value = 7497.875
$ws.Range("K137").Value = 34947.198
$ws.Range("L137").Value = 22493.625
$ws.Range("M137").Value = -32397.198
$ws.Range("N137").Value = -27593.625

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 23811772
$ws.Range("I138").Value = 35715420
$ws.Range("J138").Value = 4481.357
$ws.Range("K138").Value = 107146260
$ws.Range("L138").Value = 13444.071
$ws.Range("M138").Value = -107141120
$ws.Range("N138").Value = -23724.071

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4832.9717
$ws.Range("I32").Value = 4951.8823
$ws.Range("K32").Value = 4951.8823
$ws.Range("M32").Value = -4664.8823

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7387.9585
$ws.Range("I45").Value = 7650.1113
$ws.Range("J45").Value = 6601.5
$ws.Range("K45").Value = 7650.1113
$ws.Range("L45").Value = 6601.5
$ws.Range("M45").Value = -7273.1113
$ws.Range("N45").Value = -7355.5

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3674.7144
$ws.Range("I110").Value = 4274.9414
$ws.Range("J110").Value = 1123.75
$ws.Range("K110").Value = 4274.9414
$ws.Range("L110").Value = 1123.75
$ws.Range("M110").Value = -2229.9414
$ws.Range("N110").Value = -5213.75

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1689.4667
$ws.Range("I99").Value = 1525.8846
$ws.Range("J99").Value = 2752.75
$ws.Range("K99").Value = 1525.8846
$ws.Range("L99").Value = 2752.75
$ws.Range("M99").Value = -27.88460000000009
$ws.Range("N99").Value = -5748.75

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 38.809525
$ws.Range("I7").Value = 65.818184
$ws.Range("J7").Value = 9.1
$ws.Range("K7").Value = 65.818184
$ws.Range("L7").Value = 9.1
$ws.Range("M7").Value = 47.181816
$ws.Range("N7").Value = -235.1

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2274.8438
$ws.Range("I31").Value = 1833.1305
$ws.Range("J31").Value = 3403.6667
$ws.Range("K31").Value = 1833.1305
$ws.Range("L31").Value = 3403.6667
$ws.Range("M31").Value = -1538.1305
$ws.Range("N31").Value = -3993.6667

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2274.8438
$ws.Range("I34").Value = 1833.1305
$ws.Range("J34").Value = 3403.6667
$ws.Range("K34").Value = 1833.1305
$ws.Range("L34").Value = 3403.6667
$ws.Range("M34").Value = -1631.1305
$ws.Range("N34").Value = -3807.6667

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7171.9287
$ws.Range("I132").Value = 3580.0908
$ws.Range("J132").Value = 20342
$ws.Range("K132").Value = 10740.2724
$ws.Range("L132").Value = 61026
$ws.Range("M132").Value = -8210.2724
$ws.Range("N132").Value = -66086

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2649.54
$ws.Range("I134").Value = 2488.761
$ws.Range("K134").Value = 7466.282999999999
$ws.Range("M134").Value = -4931.282999999999

# CUL row 6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 257.66666
$ws.Range("I6").Value = 127.875
$ws.Range("J6").Value = 517.25
$ws.Range("K6").Value = 383.625
$ws.Range("L6").Value = 1551.75
$ws.Range("M6").Value = -270.625
$ws.Range("N6").Value = -1777.75

# CUL row 37
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 74419
$ws.Range("J37").Value = 74419
$ws.Range("L37").Value = 223257
$ws.Range("N37").Value = -223481

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1085.3846
$ws.Range("I113").Value = 875.1875
$ws.Range("J113").Value = 1421.7
$ws.Range("K113").Value = 2625.5625
$ws.Range("L113").Value = 4265.1
$ws.Range("M113").Value = -455.5625
$ws.Range("N113").Value = -8605.1

# GSM row 11
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3010000
$ws.Range("I11").Value = 3010000
$ws.Range("K11").Value = 3010000
$ws.Range("M11").Value = -3009861

# GSM row 14
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 1889149.2
$ws.Range("I14").Value = 3000639.2
$ws.Range("J14").Value = 36666
$ws.Range("K14").Value = 3000639.2
$ws.Range("L14").Value = 36666
$ws.Range("M14").Value = -3000471.2
$ws.Range("N14").Value = -37002

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18620.285
$ws.Range("I70").Value = 20890.334
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 20890.334
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -20620.334
$ws.Range("N70").Value = -5540

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 18620.285
$ws.Range("I73").Value = 20890.334
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 20890.334
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -19954.334
$ws.Range("N73").Value = -6872

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2308.8667
$ws.Range("I126").Value = 1864.0769
$ws.Range("J126").Value = 5200
$ws.Range("K126").Value = 5592.2307
$ws.Range("L126").Value = 15600
$ws.Range("M126").Value = -3122.2307
$ws.Range("N126").Value = -20540

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1207.9166
$ws.Range("I40").Value = 1207.9166
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1207.9166
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1071.9166
$ws.Range("N40").ClearContents()

# LTW row 106
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 23497.5
$ws.Range("J106").Value = 23497.5
$ws.Range("L106").Value = 23497.5
$ws.Range("N106").Value = -26021.5

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1246.7273
$ws.Range("I107").Value = 1351.6
$ws.Range("J107").Value = 198
$ws.Range("K107").Value = 4054.8
$ws.Range("L107").Value = 594
$ws.Range("M107").Value = -2134.8
$ws.Range("N107").Value = -4434

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3992.0344
$ws.Range("I132").Value = 3011.9583
$ws.Range("J132").Value = 8696.4
$ws.Range("K132").Value = 9035.874899999999
$ws.Range("L132").Value = 26089.2
$ws.Range("M132").Value = -6505.874899999999
$ws.Range("N132").Value = -31149.2

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2002.0222
$ws.Range("I136").Value = 2081.2327
$ws.Range("J136").Value = 299
$ws.Range("K136").Value = 6243.6981
$ws.Range("L136").Value = 897
$ws.Range("M136").Value = -3693.6981
$ws.Range("N136").Value = -5997
